$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 63

# Row 3
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = 0.1363636363636364
$ws.Range("E3").Value = 0.36
$ws.Range("F3").Value = 0.1978021978021978

# Row 4
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = 0.06451612903225806
$ws.Range("E4").Value = 0.2307692307692308
$ws.Range("F4").Value = 0.1008403361344538

# Row 5
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 0.1063829787234043
$ws.Range("E5").Value = 0.1666666666666667
$ws.Range("F5").Value = 0.1298701298701299

# Row 6
$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.05357142857142857
$ws.Range("F6").Value = 0.1016949152542373

# Row 7
$ws.Range("C7").Value = 9
